# Update the worksheet data per the revised "Cessazione unione civile e convivenze" table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the shared-string text / values for rows 6-9, and delete old row 10.
$ws.Cells.Item(6, 1).Value = 31
$ws.Cells.Item(6, 2).Value = "Accordo tra le parti"

$ws.Cells.Item(7, 1).Value = 32
$ws.Cells.Item(7, 2).Value = "Recesso unilaterale"

$ws.Cells.Item(8, 1).Value = 34
$ws.Cells.Item(8, 2).Value = "Matrimonio/unione civile"

$ws.Cells.Item(9, 1).Value = 98
$ws.Cells.Item(9, 2).Value = "Decesso del convivente/unito civilmente"

# Remove the now-obsolete 10th row entirely.
$ws.Rows.Item(10).Delete()

# Widen column B to fit the new, longer text (best-fit equivalent).
$ws.Columns.Item(2).ColumnWidth = 41.33

# Move the active selection to reflect the saved view state.
$ws.Range("B8").Select()

# Restore the page setup (paper size / orientation) recorded for this sheet.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

$wb.Save()
